# Insert a new data row at row 416 (pushes existing rows 416:496 down to
# 417:497) and populate it with the new "Ají" price record for this week.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(416).Insert()

$ws.Cells.Item(416, 1).Value2  = 10
$ws.Cells.Item(416, 2).Value2  = 'Vega Modelo de Temuco'
$ws.Cells.Item(416, 3).Value2  = 'La Araucanía'
$ws.Cells.Item(416, 4).Value2  = 44504
$ws.Cells.Item(416, 5).Value2  = 9
$ws.Cells.Item(416, 6).Value2  = 100112021
$ws.Cells.Item(416, 7).Value2  = 'Ají'
$ws.Cells.Item(416, 8).Value2  = 'Inferno'
$ws.Cells.Item(416, 9).Value2  = 'Primera'
$ws.Cells.Item(416, 10).Value2 = 65
$ws.Cells.Item(416, 11).Value2 = 34000
$ws.Cells.Item(416, 12).Value2 = 34000
$ws.Cells.Item(416, 13).Value2 = 34000
$ws.Cells.Item(416, 14).Value2 = '$/caja 15 kilos'
$ws.Cells.Item(416, 15).Value2 = 'Región de Arica y Parinacota'
$ws.Cells.Item(416, 16).Value2 = 2267
$ws.Cells.Item(416, 17).Value2 = 15
$ws.Cells.Item(416, 18).Value2 = 'Hortaliza'
